# Applies the four textual/formatting changes described in the commit:
#  1) "session is started in iLab..." sentence: bold "customer's identity "
#     and bold "time the session h".
#  2) "In iLab the needed configuration..." paragraph: underline
#     " needed configuration is stored in an Interlock object".
#  3) "For each possible combination..." paragraph: underline
#     " to be scheduled and ran together".
#  4) "The iLab Interlock module..." paragraph: underline the whole
#     sentence up through "...Internet at large".

$d = $word.ActiveDocument

function Mark-SubRange($container, $substr, $mode) {
    # Work off a duplicate of $container so repeated Find calls don't
    # disturb the caller's range, then apply direct character formatting
    # to just the matched sub-range (Find-derived ranges apply formatting
    # narrowly; ranges built from raw offsets via $d.Range() do not).
    $r = $container.Duplicate
    $ok = $r.Find.Execute($substr)
    if (-not $ok) {
        throw "substring not found: [$substr]"
    }
    if ($mode -eq "bold") {
        $r.Font.Bold = $true
    } elseif ($mode -eq "underline") {
        $r.Font.Underline = 1
    }
    return $r
}

# ---------------------------------------------------------------------
# 1) customer/time bolding
# ---------------------------------------------------------------------
$sentence1 = $d.Content
$found1 = $sentence1.Find.Execute("session is started in iLab, we will record the customer’s identity and the time the session has started as part of the reservation. ")
if (-not $found1) { throw "hunk1 anchor sentence not found" }
$anchor1 = $sentence1.Duplicate

$null = Mark-SubRange $anchor1 "customer’s identity " "bold"
$null = Mark-SubRange $anchor1 "time the session h" "bold"

# ---------------------------------------------------------------------
# 2) "In iLab the ... Interlock object" underline
# ---------------------------------------------------------------------
$sentence2 = $d.Content
$found2 = $sentence2.Find.Execute("In iLab the needed configuration is stored in an Interlock object, which can have multiple channels. Each channel can be associated with a s")
if (-not $found2) { throw "hunk2 anchor sentence not found" }
$anchor2 = $sentence2.Duplicate

$null = Mark-SubRange $anchor2 " needed configuration is stored in an Interlock object" "underline"

# ---------------------------------------------------------------------
# 3) "For each possible combination ... ran together" underline
# ---------------------------------------------------------------------
$sentence3 = $d.Content
$found3 = $sentence3.Find.Execute("For each possible combination of component and lab that will need to be scheduled and ran together, we need to configure ")
if (-not $found3) { throw "hunk3 anchor sentence not found" }
$anchor3 = $sentence3.Duplicate

$null = Mark-SubRange $anchor3 " to be scheduled and ran together" "underline"

# ---------------------------------------------------------------------
# 4) "The iLab Interlock module ... at large" underline
# ---------------------------------------------------------------------
$sentence4 = $d.Content
$found4 = $sentence4.Find.Execute("The iLab Interlock module does not require for the API endpoints of the interlock devices to be exposed to the Internet at large. This is the purpose of the iLab Bridge")
if (-not $found4) { throw "hunk4 anchor sentence not found" }
$anchor4 = $sentence4.Duplicate

$null = Mark-SubRange $anchor4 "The iLab Interlock module does not require for the API endpoints of the interlock devices to be exposed to the Internet at large" "underline"

Write-Output "All four edits applied successfully."
